$wb = $excel.ActiveWorkbook

# Rename the "Include from Tempcodes" sheet to "Include from CareSocialCodes"
$metaSheet = $wb.Worksheets.Item("Metadata")
$codesSheet = $wb.Worksheets.Item("Include from Tempcodes")
$codesSheet.Name = "Include from CareSocialCodes"

# Update Metadata sheet values
$metaSheet.Range("B3").Value = "2.0.0"
$metaSheet.Range("B8").Value = "2024-06-14T10:48:54+02:00"
$metaSheet.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# Update the renamed sheet's System URI value
$codesSheet.Range("B8").Value = "http://fhir.kl.dk/term/CodeSystem/CareSocialCodes"
